# "Generate Report for Handoff" - regenerate handoff identifiers/timestamps
# for the localization-status report.
#
# The old handoff run used commit-id "fa88caae-6336-4160-80e5-73581798e10a"
# and xlf-hash "76747d810813924f630d7e92e6e8dee07573ca57"; this run produced
# "81681032-9eb1-49b4-80ad-a0ca5d3a8a13" / "df4a2450f69c9e296db9568317ef0d07706dc9ab"
# plus refreshed handoff timestamps. Update both the cell text and the
# matching hyperlink display text (the hyperlink target URLs themselves are
# not touched by this regeneration).

$wb = $excel.ActiveWorkbook

$oldId  = "fa88caae-6336-4160-80e5-73581798e10a"
$newId  = "81681032-9eb1-49b4-80ad-a0ca5d3a8a13"
$oldHash = "76747d810813924f630d7e92e6e8dee07573ca57"
$newHash = "df4a2450f69c9e296db9568317ef0d07706dc9ab"

$oldMd  = "$oldId.md"
$newMd  = "$newId.md"

# Per-sheet cell edits: sheet name -> list of (cell, oldValue, newValue)
$edits = @{
  "Overview" = @(
    @{ Cell = "A2"; Old = $oldMd; New = $newMd },
    @{ Cell = "D2"; Old = "2016-46-20 00:46:04"; New = "2016-46-20 00:46:31" }
  )
  "zh-cn" = @(
    @{ Cell = "A2"; Old = $oldMd; New = $newMd },
    @{ Cell = "D2"; Old = "$oldId.$oldHash.zh-cn.xlf"; New = "$newId.$newHash.zh-cn.xlf" },
    @{ Cell = "E2"; Old = "2016-03-20 00:46:00"; New = "2016-03-20 00:46:28" }
  )
  "de-de" = @(
    @{ Cell = "A2"; Old = $oldMd; New = $newMd },
    @{ Cell = "D2"; Old = "$oldId.$oldHash.de-de.xlf"; New = "$newId.$newHash.de-de.xlf" },
    @{ Cell = "E2"; Old = "2016-03-20 00:46:04"; New = "2016-03-20 00:46:31" }
  )
}

foreach ($sheetName in $edits.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Update the cell values that carry the old identifiers/timestamps.
  foreach ($edit in $edits[$sheetName]) {
    $ws.Range($edit.Cell).Value = $edit.New
  }

  # Keep hyperlink display text (shown label) in sync with the new cell
  # text wherever it still references one of the old strings. Enumerating
  # the collection (rather than indexing with .Item) updates the existing
  # <hyperlink> entries in place instead of appending new ones.
  foreach ($hl in $ws.Hyperlinks) {
    $display = $hl.TextToDisplay
    foreach ($edit in $edits[$sheetName]) {
      if ($display -eq $edit.Old) {
        $hl.TextToDisplay = $edit.New
        break
      }
    }
  }
}
